# billing/in/rates.xlsx — "added ui and few changes for this"
#
# The rates sheet gained one more product row ("baba") underneath the
# existing Navel/Blood/Mandarin/... rows, with a Rate of 95 and a Scope
# value of 45 (numeric, like the other recently-added Scope values in
# rows 10-13). After typing the row the user's selection ended up one
# row below/at column C (C15), which is what Excel leaves selected after
# finishing data entry on row 14's last populated column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14: Product / Rate / Scope
$ws.Range("A14").Value = "baba"
$ws.Range("B14").Value = 95
$ws.Range("C14").Value = 45

# Leave the selection where the user's cursor would land after entering
# the new row of data.
[void]$ws.Range("C15").Select()
